$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the pairing matrix (minutes paired together) for Team 01.
$ws.Range("F15").Value = 240
$ws.Range("E16").Value = 480
$ws.Range("F16").Value = 120
$ws.Range("G17").Value = 240
$ws.Range("H19").Value = 480
$ws.Range("G20").Value = 480

# Move the active selection as recorded when the workbook was last saved.
$ws.Range("I18").Select()
